$wb = $excel.ActiveWorkbook

# Template sheet to clone for the new country
$spain = $wb.Worksheets.Item("Spain")

# Reset the selection on the Spain sheet to the full used range (no specific
# active cell) before it loses focus - mirrors the author's saved state.
$spain.Range("A1:D11").Select()

# Copy "Spain" and place the copy right after it
$spain.Copy([System.Reflection.Missing]::Value, $spain)

# The newly created sheet is now active; rename it
$turkey = $wb.ActiveSheet
$turkey.Name = "Turkey"

# Fill in the country-specific values
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3307"

# The copied sheet carried over Spain's cached (non-pinned) row heights;
# re-fit them so unused wrap-height goes back to the sheet default.
$turkey.Rows.Item(3).RowHeight = 14.4
$turkey.Rows.Item(3).AutoFit()
$turkey.Rows.Item(4).RowHeight = 14.4
$turkey.Rows.Item(4).AutoFit()
$turkey.Rows.Item(5).RowHeight = 14.4
$turkey.Rows.Item(5).AutoFit()

# Match the saved selection on the new sheet
$turkey.Range("G15").Select()
